{"js": "// Update the worksheet date and all 25 three-digit-by-one-digit\n// multiplication problems to the new values, per the commit diff.\nconst replacements = [\n  [\"2026-01-12 Monday\", \"2026-01-13 Tuesday\"],\n  [\"327\u00d79=2943\", \"431\u00d78=3448\"],\n  [\"356\u00d77=2492\", \"905\u00d78=7240\"],\n  [\"305\u00d78=2440\", \"528\u00d72=1056\"],\n  [\"367\u00d74=1468\", \"795\u00d78=6360\"],\n  [\"568\u00d73=1704\", \"384\u00d78=3072\"],\n  [\"810\u00d73=2430\", \"847\u00d72=1694\"],\n  [\"469\u00d76=2814\", \"718\u00d72=1436\"],\n  [\"856\u00d73=2568\", \"599\u00d79=5391\"],\n  [\"956\u00d76=5736\", \"515\u00d78=4120\"],\n  [\"881\u00d73=2643\", \"477\u00d72=954\"],\n  [\"549\u00d75=2745\", \"670\u00d74=2680\"],\n  [\"714\u00d78=5712\", \"970\u00d75=4850\"],\n  [\"256\u00d74=1024\", \"868\u00d75=4340\"],\n  [\"579\u00d75=2895\", \"738\u00d72=1476\"],\n  [\"986\u00d75=4930\", \"357\u00d76=2142\"],\n  [\"370\u00d78=2960\", \"922\u00d79=8298\"],\n  [\"932\u00d76=5592\", \"650\u00d76=3900\"],\n  [\"824\u00d77=5768\", \"362\u00d76=2172\"],\n  [\"978\u00d76=5868\", \"757\u00d78=6056\"],\n  [\"933\u00d75=4665\", \"475\u00d78=3800\"],\n  [\"967\u00d77=6769\", \"818\u00d73=2454\"],\n  [\"110\u00d77=770\", \"611\u00d73=1833\"],\n  [\"754\u00d73=2262\", \"559\u00d72=1118\"],\n  [\"569\u00d72=1138\", \"800\u00d79=7200\"],\n  [\"988\u00d75=4940\", \"485\u00d78=3880\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 three-digit-by-one-digit\n# multiplication problems to the new values, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-12 Monday\", \"2026-01-13 Tuesday\"),\n    @(\"327\u00d79=2943\", \"431\u00d78=3448\"),\n    @(\"356\u00d77=2492\", \"905\u00d78=7240\"),\n    @(\"305\u00d78=2440\", \"528\u00d72=1056\"),\n    @(\"367\u00d74=1468\", \"795\u00d78=6360\"),\n    @(\"568\u00d73=1704\", \"384\u00d78=3072\"),\n    @(\"810\u00d73=2430\", \"847\u00d72=1694\"),\n    @(\"469\u00d76=2814\", \"718\u00d72=1436\"),\n    @(\"856\u00d73=2568\", \"599\u00d79=5391\"),\n    @(\"956\u00d76=5736\", \"515\u00d78=4120\"),\n    @(\"881\u00d73=2643\", \"477\u00d72=954\"),\n    @(\"549\u00d75=2745\", \"670\u00d74=2680\"),\n    @(\"714\u00d78=5712\", \"970\u00d75=4850\"),\n    @(\"256\u00d74=1024\", \"868\u00d75=4340\"),\n    @(\"579\u00d75=2895\", \"738\u00d72=1476\"),\n    @(\"986\u00d75=4930\", \"357\u00d76=2142\"),\n    @(\"370\u00d78=2960\", \"922\u00d79=8298\"),\n    @(\"932\u00d76=5592\", \"650\u00d76=3900\"),\n    @(\"824\u00d77=5768\", \"362\u00d76=2172\"),\n    @(\"978\u00d76=5868\", \"757\u00d78=6056\"),\n    @(\"933\u00d75=4665\", \"475\u00d78=3800\"),\n    @(\"967\u00d77=6769\", \"818\u00d73=2454\"),\n    @(\"110\u00d77=770\", \"611\u00d73=1833\"),\n    @(\"754\u00d73=2262\", \"559\u00d72=1118\"),\n    @(\"569\u00d72=1138\", \"800\u00d79=7200\"),\n    @(\"988\u00d75=4940\", \"485\u00d78=3880\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
